# The document originally holds "Create a VPC with the specified CIDR blocks"
# and ": " as two separate runs (same run formatting). The edit merges them
# into a single run/text: "Create a VPC with the specified CIDR blocks: ".
$d = $word.ActiveDocument

$find = $d.Content.Find
$find.ClearFormatting()
$find.Replacement.ClearFormatting()
$find.Execute("Create a VPC with the specified CIDR blocks: ", $true, $false, $false, $false, $false,
              $true, 1, $false, "Create a VPC with the specified CIDR blocks: ", 2)
